{"js": "// Update the division problems in the table to new values, per the diff.\nconst replacements = [\n  [\"979\u00f75=\", \"806\u00f78=\"],\n  [\"256\u00f76=\", \"114\u00f75=\"],\n  [\"470\u00f74=\", \"898\u00f74=\"],\n  [\"473\u00f74=\", \"892\u00f72=\"],\n  [\"301\u00f77=\", \"929\u00f72=\"],\n  [\"795\u00f73=\", \"739\u00f74=\"],\n  [\"802\u00f74=\", \"939\u00f78=\"],\n  [\"810\u00f79=\", \"883\u00f78=\"],\n  [\"286\u00f77=\", \"346\u00f72=\"],\n  [\"715\u00f78=\", \"246\u00f75=\"],\n  [\"512\u00f73=\", \"925\u00f78=\"],\n  [\"345\u00f79=\", \"584\u00f72=\"],\n  [\"516\u00f72=\", \"641\u00f75=\"],\n  [\"218\u00f76=\", \"110\u00f78=\"],\n  [\"164\u00f77=\", \"924\u00f77=\"],\n  [\"226\u00f77=\", \"621\u00f77=\"],\n  [\"449\u00f76=\", \"196\u00f79=\"],\n  [\"314\u00f78=\", \"577\u00f79=\"],\n  [\"825\u00f75=\", \"496\u00f77=\"],\n  [\"236\u00f77=\", \"834\u00f74=\"],\n  [\"284\u00f79=\", \"878\u00f75=\"],\n  [\"509\u00f79=\", \"558\u00f73=\"],\n  [\"967\u00f73=\", \"490\u00f78=\"],\n  [\"627\u00f75=\", \"855\u00f76=\"],\n  [\"345\u00f74=\", \"956\u00f75=\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const searchResults = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  searchResults.load(\"items\");\n  await context.sync();\n\n  for (const result of searchResults.items) {\n    result.insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# Update the division problems in the table to new values, per the diff.\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"979\u00f75=\", \"806\u00f78=\"),\n    @(\"256\u00f76=\", \"114\u00f75=\"),\n    @(\"470\u00f74=\", \"898\u00f74=\"),\n    @(\"473\u00f74=\", \"892\u00f72=\"),\n    @(\"301\u00f77=\", \"929\u00f72=\"),\n    @(\"795\u00f73=\", \"739\u00f74=\"),\n    @(\"802\u00f74=\", \"939\u00f78=\"),\n    @(\"810\u00f79=\", \"883\u00f78=\"),\n    @(\"286\u00f77=\", \"346\u00f72=\"),\n    @(\"715\u00f78=\", \"246\u00f75=\"),\n    @(\"512\u00f73=\", \"925\u00f78=\"),\n    @(\"345\u00f79=\", \"584\u00f72=\"),\n    @(\"516\u00f72=\", \"641\u00f75=\"),\n    @(\"218\u00f76=\", \"110\u00f78=\"),\n    @(\"164\u00f77=\", \"924\u00f77=\"),\n    @(\"226\u00f77=\", \"621\u00f77=\"),\n    @(\"449\u00f76=\", \"196\u00f79=\"),\n    @(\"314\u00f78=\", \"577\u00f79=\"),\n    @(\"825\u00f75=\", \"496\u00f77=\"),\n    @(\"236\u00f77=\", \"834\u00f74=\"),\n    @(\"284\u00f79=\", \"878\u00f75=\"),\n    @(\"509\u00f79=\", \"558\u00f73=\"),\n    @(\"967\u00f73=\", \"490\u00f78=\"),\n    @(\"627\u00f75=\", \"855\u00f76=\"),\n    @(\"345\u00f74=\", \"956\u00f75=\")\n)\n\nforeach ($pair in $replacements) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Text = $oldText\n    $find.Replacement.ClearFormatting()\n    $find.Replacement.Text = $newText\n    $find.Execute([ref]$oldText, [ref]$true, [ref]$false, [ref]$false, [ref]$false, [ref]$false, [ref]$true, [ref]1, [ref]$true, [ref]$newText, [ref]2) | Out-Null\n}\n"}
